$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.613.18'
$ws.Range("E2").Value = '  -1.12%  '

$ws.Range("D3").Value = '3.522.33'
$ws.Range("E3").Value = '  -3.33%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '201.72'
$ws.Range("E5").Value = '  +2.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '551.79'
$ws.Range("E6").Value = '  -4.03%  '

$ws.Range("D7").Value = '3.512.79'
$ws.Range("E7").Value = '  -3.41%  '

$ws.Range("E8").Value = '  -1.35%  '

$ws.Range("E9").Value = '  +0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.653'
$ws.Range("E10").Value = '  -3.78%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '62.70'
$ws.Range("E11").Value = '  +10.84%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.143'
$ws.Range("E12").Value = '  -7.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000270'
$ws.Range("E13").Value = '  -9.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.86'
$ws.Range("E14").Value = '  -2.57%  '

$ws.Range("D15").Value = '4.086.29'
$ws.Range("E15").Value = '  -3.18%  '

$ws.Range("D16").Value = '3.525.36'
$ws.Range("E16").Value = '  -3.09%  '

$ws.Range("E17").Value = '  -1.78%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.44'
$ws.Range("E18").Value = '  -0.98%  '

$ws.Range("D19").Value = '67.410.90'
$ws.Range("E19").Value = '  -1.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.85'
$ws.Range("E20").Value = '  -5.83%  '

$ws.Range("E21").Value = '  -5.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '393.44'
$ws.Range("E22").Value = '  -2.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.99'
$ws.Range("E23").Value = '  -9.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.99'
$ws.Range("E24").Value = '  -6.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.99'
$ws.Range("E25").Value = '  -2.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.90'
$ws.Range("E26").Value = '  +0.10%  '

$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.83'
$ws.Range("E27").Value = '  -4.91%  '

$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.26'
$ws.Range("E28").Value = '  -3.41%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.85'
$ws.Range("E29").Value = '  -4.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '717.93'
$ws.Range("E30").Value = '  +3.96%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.12'
$ws.Range("E31").Value = '  -2.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.08'
$ws.Range("E32").Value = '  -14.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.74'
$ws.Range("E33").Value = '  -4.45%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.77'
$ws.Range("E34").Value = '  -1.63%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.111'
$ws.Range("E35").Value = '  -5.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.56'
$ws.Range("E36").Value = '  -10.00%  '

$ws.Range("E37").Value = '  +0.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.396'
$ws.Range("E38").Value = '  -8.24%  '

$ws.Range("B39").Value = 'ThetaToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.03'
$ws.Range("E39").Value = '  -4.06%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.131'
$ws.Range("E40").Value = '  -6.91%  '

$ws.Range("D41").Value = '3.082.81'
$ws.Range("E41").Value = '  -4.30%  '

$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").Value = '0.0₃0681'
$ws.Range("E43").Value = '  -14.62%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.82'
$ws.Range("E44").Value = '  +7.24%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.54'
$ws.Range("E45").Value = '  -14.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0409'
$ws.Range("E46").Value = '  -2.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.68'
$ws.Range("E47").Value = '  -12.33%  '

$ws.Range("E48").Value = '  -3.82%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.23'
$ws.Range("E49").Value = '  -3.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.27'
$ws.Range("E50").Value = '  -7.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.88'
$ws.Range("E51").Value = '  -6.60%  '
